# Applies cached LeveProfits market-price figures updated by the scheduled Shiva price-scrape runner.
# Values correspond 1:1 to the authoritative XLSX diff; all cells are plain numeric (no formulas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 16
$ws.Range("H16").Value = 4900
$ws.Range("I16").Value = 4900
$ws.Range("K16").Value = 4900
$ws.Range("M16").Value = -4670

# ALC row 98
$ws.Range("H98").Value = 933.6667
$ws.Range("J98").Value = 3999.5
$ws.Range("L98").Value = 3999.5
$ws.Range("N98").Value = -6995.5

# ALC row 106
$ws.Range("H106").Value = 120167.89
$ws.Range("I106").Value = 172587.67
$ws.Range("K106").Value = 172587.67
$ws.Range("M106").Value = -171956.67

# ALC row 122
$ws.Range("H122").Value = 933.6667
$ws.Range("J122").Value = 3999.5
$ws.Range("L122").Value = 11998.5
$ws.Range("N122").Value = -16898.5

# ALC row 132
$ws.Range("H132").Value = 5330.115
$ws.Range("I132").Value = 2783.3157
$ws.Range("K132").Value = 8349.947100000001
$ws.Range("M132").Value = -5819.947100000001

# ALC row 138
$ws.Range("H138").Value = 1841.64
$ws.Range("I138").Value = 1275.1111
$ws.Range("J138").Value = 2160.3125
$ws.Range("K138").Value = 3825.3333
$ws.Range("L138").Value = 6480.9375
$ws.Range("M138").Value = 1314.6667
$ws.Range("N138").Value = -16760.9375

# ALC row 141
$ws.Range("H141").Value = 8889.048000000001
$ws.Range("I141").Value = 9874.615
$ws.Range("J141").Value = 7287.5
$ws.Range("K141").Value = 29623.845
$ws.Range("L141").Value = 21862.5
$ws.Range("M141").Value = -24443.845
$ws.Range("N141").Value = -32222.5

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 5508.7114
$ws.Range("I32").Value = 5247.1626
$ws.Range("J32").Value = 6758.3335
$ws.Range("K32").Value = 5247.1626
$ws.Range("L32").Value = 6758.3335
$ws.Range("M32").Value = -4960.1626
$ws.Range("N32").Value = -7332.3335

# ARM row 80
$ws.Range("H80").Value = 49999
$ws.Range("J80").Value = 49999
$ws.Range("L80").Value = 49999
$ws.Range("N80").Value = -51995

# ARM row 83
$ws.Range("H83").Value = 49999
$ws.Range("J83").Value = 49999
$ws.Range("L83").Value = 149997
$ws.Range("N83").Value = -159981

# ARM row 132
$ws.Range("H132").Value = 933.0714
$ws.Range("I132").Value = 938
$ws.Range("J132").Value = 903.5
$ws.Range("K132").Value = 2814
$ws.Range("L132").Value = 2710.5
$ws.Range("M132").Value = -284
$ws.Range("N132").Value = -7770.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 1821.3654
$ws.Range("J20").Value = 1678.9474
$ws.Range("L20").Value = 1678.9474
$ws.Range("N20").Value = -2172.9474

# BSM row 86
$ws.Range("H86").Value = 2136.925
$ws.Range("I86").Value = 1972.2
$ws.Range("J86").Value = 2411.4666
$ws.Range("K86").Value = 1972.2
$ws.Range("L86").Value = 2411.4666
$ws.Range("M86").Value = -849.2
$ws.Range("N86").Value = -4657.4666

# BSM row 89
$ws.Range("H89").Value = 2136.925
$ws.Range("I89").Value = 1972.2
$ws.Range("J89").Value = 2411.4666
$ws.Range("K89").Value = 9861
$ws.Range("L89").Value = 12057.333
$ws.Range("M89").Value = -4245
$ws.Range("N89").Value = -23289.333

# BSM row 105
$ws.Range("H105").Value = 3226.2327
$ws.Range("I105").Value = 2371.611
$ws.Range("K105").Value = 2371.611
$ws.Range("M105").Value = -624.6109999999999

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 1945.4839
$ws.Range("I31").Value = 1369.6666
$ws.Range("J31").Value = 2485.3125
$ws.Range("K31").Value = 1369.6666
$ws.Range("L31").Value = 2485.3125
$ws.Range("M31").Value = -1074.6666
$ws.Range("N31").Value = -3075.3125

# CRP row 34
$ws.Range("H34").Value = 1945.4839
$ws.Range("I34").Value = 1369.6666
$ws.Range("J34").Value = 2485.3125
$ws.Range("K34").Value = 1369.6666
$ws.Range("L34").Value = 2485.3125
$ws.Range("M34").Value = -1167.6666
$ws.Range("N34").Value = -2889.3125

# CRP row 132
$ws.Range("H132").Value = 5007.4375
$ws.Range("I132").Value = 2114.5476
$ws.Range("K132").Value = 6343.6428
$ws.Range("M132").Value = -3813.6428

$ws = $wb.Worksheets.Item("CUL")
# CUL row 12
$ws.Range("H12").Value = 118.32
$ws.Range("I12").Value = 91.59999999999999
$ws.Range("J12").Value = 125
$ws.Range("K12").Value = 274.8
$ws.Range("L12").Value = 375
$ws.Range("M12").Value = -101.8
$ws.Range("N12").Value = -721

# CUL row 14
$ws.Range("H14").Value = 421.93332
$ws.Range("I14").Value = 421.93332
$ws.Range("K14").Value = 1265.79996
$ws.Range("M14").Value = -1092.79996

# CUL row 33
$ws.Range("H33").Value = 365
$ws.Range("I33").Value = 94
$ws.Range("K33").Value = 564
$ws.Range("M33").Value = -281

# CUL row 122
$ws.Range("H122").Value = 1071.8889
$ws.Range("I122").Value = 672.2857
$ws.Range("J122").Value = 1326.1818
$ws.Range("K122").Value = 6050.571300000001
$ws.Range("L122").Value = 11935.6362
$ws.Range("M122").Value = -3600.571300000001
$ws.Range("N122").Value = -16835.6362

$ws = $wb.Worksheets.Item("GSM")
# GSM row 43
$ws.Range("H43").Value = 13833.333
$ws.Range("I43").Value = 11600
$ws.Range("K43").Value = 11600
$ws.Range("M43").Value = -11449

# GSM row 80
$ws.Range("H80").Value = 24214076
$ws.Range("I80").Value = 30670052
$ws.Range("K80").Value = 30670052
$ws.Range("M80").Value = -30669054

# GSM row 83
$ws.Range("H83").Value = 24214076
$ws.Range("I83").Value = 30670052
$ws.Range("K83").Value = 153350260
$ws.Range("M83").Value = -153345268

# GSM row 97
$ws.Range("H97").Value = 450.36
$ws.Range("I97").Value = 362.81818
$ws.Range("J97").Value = 1092.3334
$ws.Range("K97").Value = 362.81818
$ws.Range("L97").Value = 1092.3334
$ws.Range("M97").Value = 133.18182
$ws.Range("N97").Value = -2084.3334

# GSM row 132
$ws.Range("H132").Value = 3549
$ws.Range("I132").Value = 3610.7778
$ws.Range("J132").Value = 2993
$ws.Range("K132").Value = 10832.3334
$ws.Range("L132").Value = 8979
$ws.Range("M132").Value = -8302.3334
$ws.Range("N132").Value = -14039

$ws = $wb.Worksheets.Item("LTW")
# LTW row 93
$ws.Range("H93").Value = 1346.2
$ws.Range("I93").Value = 1109.4166
$ws.Range("K93").Value = 1109.4166
$ws.Range("M93").Value = 138.5834

# LTW row 122
$ws.Range("H122").Value = 3313.6
$ws.Range("I122").Value = 3429.8
$ws.Range("J122").Value = 2848.8
$ws.Range("K122").Value = 10289.4
$ws.Range("L122").Value = 8546.400000000001
$ws.Range("M122").Value = -7839.400000000001
$ws.Range("N122").Value = -13446.4

# LTW row 132
$ws.Range("H132").Value = 52352.7
$ws.Range("I132").Value = 61062.117
$ws.Range("K132").Value = 183186.351
$ws.Range("M132").Value = -180656.351

# LTW row 134
$ws.Range("H134").Value = 124331
$ws.Range("J134").Value = 124331
$ws.Range("L134").Value = 124331
$ws.Range("N134").Value = -134471

# LTW row 136
$ws.Range("H136").Value = 10264.286
$ws.Range("I136").Value = 10760.462
$ws.Range("K136").Value = 32281.386
$ws.Range("M136").Value = -29731.386

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 3429.7
$ws.Range("I122").Value = 3429.7
$ws.Range("K122").Value = 10289.1
$ws.Range("M122").Value = -7839.099999999999

Write-Host "Applied all profit sheet updates"
